# Apply the commit's substantive changes:
#  - fixed bug that caused scalar data from local file not being loaded
#    in certain situations: the "scalars" sheet is renamed to "_scalars"
#    (leading underscore), and the scalar values "mins" and "beta" are
#    corrected.

$wb = $excel.ActiveWorkbook

# Rename the "scalars" sheet to "_scalars"
$scalars = $wb.Worksheets.Item("scalars")
$scalars.Name = "_scalars"

# Correct the scalar values on the renamed sheet:
#   mins (row 3, col C) : 100 -> 10
#   beta (row 4, col C) : 0.95 -> 0.8
$scalars.Range("C3").Value = 10
$scalars.Range("C4").Value = 0.8

# Leave the active selection on D6, matching the saved workbook state.
[void]$scalars.Range("D6").Select()
